# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Ají" at rows 458-459 of Sheet1,
# pushing the existing rows 458:520 down to 460:522.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 458 (shifts 458:520 -> 460:522)
$ws.Rows("458:459").Insert()

# New row 458: Inferno / Primera, Región de Arica y Parinacota
$ws.Range("A458").Value = 8
$ws.Range("B458").Value = "Terminal La Palmera de La Serena"
$ws.Range("C458").Value = "Coquimbo"
$ws.Range("D458").Value = 45154
$ws.Range("E458").Value = 4
$ws.Range("F458").Value = 100112021
$ws.Range("G458").Value = "Ají"
$ws.Range("H458").Value = "Inferno"
$ws.Range("I458").Value = "Primera"
$ws.Range("J458").Value = 560
$ws.Range("K458").Value = 15000
$ws.Range("L458").Value = 16000
$ws.Range("M458").Value = 15500
$ws.Range("N458").Value = "`$/caja 10 kilos"
$ws.Range("O458").Value = "Región de Arica y Parinacota"
$ws.Range("P458").Value = 1550
$ws.Range("Q458").Value = 10
$ws.Range("R458").Value = "Hortaliza"

# New row 459: Inferno / Segunda, Región de Arica y Parinacota
$ws.Range("A459").Value = 8
$ws.Range("B459").Value = "Terminal La Palmera de La Serena"
$ws.Range("C459").Value = "Coquimbo"
$ws.Range("D459").Value = 45154
$ws.Range("E459").Value = 4
$ws.Range("F459").Value = 100112021
$ws.Range("G459").Value = "Ají"
$ws.Range("H459").Value = "Inferno"
$ws.Range("I459").Value = "Segunda"
$ws.Range("J459").Value = 320
$ws.Range("K459").Value = 10000
$ws.Range("L459").Value = 11000
$ws.Range("M459").Value = 10500
$ws.Range("N459").Value = "`$/caja 10 kilos"
$ws.Range("O459").Value = "Región de Arica y Parinacota"
$ws.Range("P459").Value = 1050
$ws.Range("Q459").Value = 10
$ws.Range("R459").Value = "Hortaliza"
